$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 6
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 4
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 2
$ws.Range("B20").Value = 1
$ws.Range("B21").Value = 11
$ws.Range("B22").Value = 4
$ws.Range("B23").Value = 5
$ws.Range("B24").Value = 5
$ws.Range("B27").Value = 4
$ws.Range("B28").Value = 4
$ws.Range("B29").Value = 2
$ws.Range("B30").Value = 2
$ws.Range("B31").Value = 5
$ws.Range("B32").Value = 2
$ws.Range("B33").Value = 2
$ws.Range("B34").Value = 5
$ws.Range("B35").Value = 4
$ws.Range("B36").Value = 6
$ws.Range("B37").Value = 2
$ws.Range("B39").Value = 3
$ws.Range("B40").Value = 2
$ws.Range("B41").Value = 4
$ws.Range("B42").Value = 8
$ws.Range("B43").Value = 8
$ws.Range("B44").Value = 10
$ws.Range("B45").Value = 6
$ws.Range("B46").Value = 4
$ws.Range("B47").Value = 4
$ws.Range("B48").Value = 1
$ws.Range("B49").Value = 1
$ws.Range("B50").Value = 6
$ws.Range("B51").Value = 5
$ws.Range("B52").Value = 14
$ws.Range("B53").Value = 10
$ws.Range("B54").Value = 6
$ws.Range("B55").Value = 2
$ws.Range("B56").Value = 0
$ws.Range("B57").Value = 0
$ws.Range("B58").Value = 2
$ws.Range("B59").Value = 3
$ws.Range("B60").Value = 3
$ws.Range("B61").Value = 1
$ws.Range("B62").Value = 5
$ws.Range("B63").Value = 5
$ws.Range("B64").Value = 0
$ws.Range("B65").Value = 1
$ws.Range("B66").Value = 5
$ws.Range("B67").Value = 3
$ws.Range("B68").Value = 2
$ws.Range("B69").Value = 3
$ws.Range("B70").Value = 3
$ws.Range("B71").Value = 1
$ws.Range("B72").Value = 1
$ws.Range("B73").Value = 4
$ws.Range("B74").Value = 4
$ws.Range("B75").Value = 2
$ws.Range("B76").Value = 2
$ws.Range("B77").Value = 5
$ws.Range("B78").Value = 2
$ws.Range("B79").Value = 2
$ws.Range("B80").Value = 2
$ws.Range("B81").Value = 0
$ws.Range("B83").Value = 1
$ws.Range("B84").Value = 3
$ws.Range("B85").Value = 0
$ws.Range("B86").Value = 6
$ws.Range("B87").Value = 2
$ws.Range("B88").Value = 1
$ws.Range("B89").Value = 2
$ws.Range("B90").Value = 1
$ws.Range("B91").Value = 1
$ws.Range("B92").Value = 0
$ws.Range("B93").Value = 1
$ws.Range("B94").Value = 0
$ws.Range("B95").Value = 5
$ws.Range("B96").Value = 1
$ws.Range("B97").Value = 0
$ws.Range("B98").Value = 1
$ws.Range("B99").Value = 5
$ws.Range("B100").Value = 3
